# Update market-price derived columns (H:N) on each job sheet.
# Values come from a refreshed market-data snapshot; row/col layout is unchanged.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 9
$ws.Range("H9").Value = 42.8
$ws.Range("I9").Value = 43.11111
$ws.Range("J9").Value = 40
$ws.Range("K9").Value = 43.11111
$ws.Range("L9").Value = 40
$ws.Range("M9").Value = 125.88889
$ws.Range("N9").Value = -378

# Row 32
$ws.Range("H32").Value = 902
$ws.Range("I32").Value = 784.2857
$ws.Range("J32").Value = 1066.8
$ws.Range("K32").Value = 784.2857
$ws.Range("L32").Value = 1066.8
$ws.Range("M32").Value = -458.2857
$ws.Range("N32").Value = -1718.8

# Row 43
$ws.Range("H43").Value = 860.5454999999999
$ws.Range("I43").Value = 980.3333
$ws.Range("J43").Value = 716.8
$ws.Range("K43").Value = 980.3333
$ws.Range("L43").Value = 716.8
$ws.Range("M43").Value = -911.3333

# Row 53
$ws.Range("H53").Value = 467.84616
$ws.Range("I53").Value = 413.57144
$ws.Range("J53").Value = 531.1667
$ws.Range("K53").Value = 413.57144
$ws.Range("L53").Value = 531.1667
$ws.Range("M53").Value = 223.42856
$ws.Range("N53").Value = -1805.1667

# Row 80
$ws.Range("H80").Value = 697.125
$ws.Range("I80").Value = 300.2857
$ws.Range("J80").Value = 1005.7778
$ws.Range("K80").Value = 900.8571000000001
$ws.Range("L80").Value = 3017.3334
$ws.Range("M80").Value = 97.14289999999994
$ws.Range("N80").Value = -5013.3334

# Row 83
$ws.Range("H83").Value = 697.125
$ws.Range("I83").Value = 300.2857
$ws.Range("J83").Value = 1005.7778
$ws.Range("K83").Value = 2702.5713
$ws.Range("L83").Value = 9052.0002
$ws.Range("M83").Value = 2289.4287
$ws.Range("N83").Value = -19036.0002

# Row 116
$ws.Range("H116").Value = 17297524
$ws.Range("I116").Value = 23062898
$ws.Range("J116").Value = 1400
$ws.Range("K116").Value = 23062898
$ws.Range("L116").Value = 1400
$ws.Range("M116").Value = -23059456
$ws.Range("N116").Value = -8284

# Row 137
$ws.Range("H137").Value = 1761.3636
$ws.Range("I137").Value = 1233.3334
$ws.Range("J137").Value = 1959.375
$ws.Range("K137").Value = 3700.0002
$ws.Range("L137").Value = 5878.125
$ws.Range("M137").Value = -1150.0002
$ws.Range("N137").Value = -10978.125

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 2148.232
$ws.Range("I32").Value = 1662.6957
$ws.Range("J32").Value = 3119.3044
$ws.Range("K32").Value = 1662.6957
$ws.Range("L32").Value = 3119.3044
$ws.Range("M32").Value = -1375.6957
$ws.Range("N32").Value = -3693.3044

# Row 36
$ws.Range("H36").Value = 30013
$ws.Range("I36").Value = 30013
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 30013
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -29667

# Row 61
$ws.Range("H61").Value = 1999.7407
$ws.Range("I61").Value = 1298.1904
$ws.Range("J61").Value = 4455.1665
$ws.Range("K61").Value = 1298.1904
$ws.Range("L61").Value = 4455.1665
$ws.Range("M61").Value = -1086.1904
$ws.Range("N61").Value = -4879.1665

# Row 136
$ws.Range("H136").Value = 1999.7407
$ws.Range("I136").Value = 1298.1904
$ws.Range("J136").Value = 4455.1665
$ws.Range("K136").Value = 3894.5712
$ws.Range("L136").Value = 13365.4995
$ws.Range("M136").Value = -1344.5712
$ws.Range("N136").Value = -18465.4995

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 580.2222
$ws.Range("I64").Value = 478.5
$ws.Range("J64").Value = 661.6
$ws.Range("K64").Value = 478.5
$ws.Range("L64").Value = 661.6
$ws.Range("M64").Value = -253.5
$ws.Range("N64").Value = -1111.6

# Row 67
$ws.Range("H67").Value = 580.2222
$ws.Range("I67").Value = 478.5
$ws.Range("J67").Value = 661.6
$ws.Range("K67").Value = 478.5
$ws.Range("L67").Value = 661.6
$ws.Range("M67").Value = 301.5
$ws.Range("N67").Value = -2221.6

# Row 86
$ws.Range("H86").Value = 13531.23
$ws.Range("I86").Value = 7766.5557
$ws.Range("J86").Value = 26501.75
$ws.Range("K86").Value = 7766.5557
$ws.Range("L86").Value = 26501.75
$ws.Range("M86").Value = -6643.5557
$ws.Range("N86").Value = -28747.75

# Row 89
$ws.Range("H89").Value = 13531.23
$ws.Range("I89").Value = 7766.5557
$ws.Range("J89").Value = 26501.75
$ws.Range("K89").Value = 38832.7785
$ws.Range("L89").Value = 132508.75
$ws.Range("M89").Value = -33216.7785
$ws.Range("N89").Value = -143740.75

# Row 107
$ws.Range("H107").Value = 2474.4
$ws.Range("I107").Value = 2505.1667
$ws.Range("J107").Value = 2428.25
$ws.Range("K107").Value = 2505.1667
$ws.Range("L107").Value = 2428.25
$ws.Range("M107").Value = -585.1667000000002
$ws.Range("N107").Value = -6268.25

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1999.2894
$ws.Range("I31").Value = 1278.1111
$ws.Range("J31").Value = 2648.35
$ws.Range("K31").Value = 1278.1111
$ws.Range("L31").Value = 2648.35
$ws.Range("M31").Value = -983.1111000000001
$ws.Range("N31").Value = -3238.35

# Row 34
$ws.Range("H34").Value = 1999.2894
$ws.Range("I34").Value = 1278.1111
$ws.Range("J34").Value = 2648.35
$ws.Range("K34").Value = 1278.1111
$ws.Range("L34").Value = 2648.35
$ws.Range("M34").Value = -1076.1111
$ws.Range("N34").Value = -3052.35

# Row 106
$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

# Row 141
$ws.Range("H141").Value = 96666.664
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 96666.664
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 96666.664
$ws.Range("N141").Value = -107026.664
$ws.Range("M141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 943.4375
$ws.Range("I23").Value = 1518.25
$ws.Range("J23").Value = 368.625
$ws.Range("K23").Value = 4554.75
$ws.Range("L23").Value = 1105.875
$ws.Range("M23").Value = -4319.75

# Row 69
$ws.Range("H69").Value = 3063.9
$ws.Range("I69").Value = 487.33334
$ws.Range("J69").Value = 3518.5881
$ws.Range("K69").Value = 1462.00002
$ws.Range("L69").Value = 10555.7643
$ws.Range("M69").Value = -651.0000199999999
$ws.Range("N69").Value = -12177.7643

# Row 72
$ws.Range("H72").Value = 3063.9
$ws.Range("I72").Value = 487.33334
$ws.Range("J72").Value = 3518.5881
$ws.Range("K72").Value = 4386.00006
$ws.Range("L72").Value = 31667.2929
$ws.Range("M72").Value = -330.0000600000003
$ws.Range("N72").Value = -39779.2929

# Row 113
$ws.Range("H113").Value = 410.26315
$ws.Range("I113").Value = 413.5
$ws.Range("J113").Value = 409.88235
$ws.Range("K113").Value = 1240.5
$ws.Range("L113").Value = 1229.64705
$ws.Range("M113").Value = 929.5
$ws.Range("N113").Value = -5569.64705

# Row 131
$ws.Range("H131").Value = 2445.506
$ws.Range("I131").Value = 378
$ws.Range("J131").Value = 2728.726
$ws.Range("K131").Value = 1134
$ws.Range("L131").Value = 8186.178
$ws.Range("M131").Value = 3906
$ws.Range("N131").Value = -18266.178

# Row 132
$ws.Range("H132").Value = 1214.1428
$ws.Range("I132").Value = 599
$ws.Range("J132").Value = 1316.6666
$ws.Range("K132").Value = 5391
$ws.Range("L132").Value = 11849.9994
$ws.Range("M132").Value = -2861
$ws.Range("N132").Value = -16909.9994

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 486.66666
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 486.66666
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 486.66666
$ws.Range("N46").Value = -862.66666

# Row 55
$ws.Range("H55").Value = 243.57143
$ws.Range("I55").Value = 75.25
$ws.Range("J55").Value = 468
$ws.Range("K55").Value = 75.25
$ws.Range("L55").Value = 468
$ws.Range("M55").Value = 97.75
$ws.Range("N55").Value = -814

# Row 122
$ws.Range("H122").Value = 3323.0417
$ws.Range("I122").Value = 2160.8
$ws.Range("J122").Value = 3628.8948
$ws.Range("K122").Value = 6482.400000000001
$ws.Range("L122").Value = 10886.6844
$ws.Range("M122").Value = -4032.400000000001
$ws.Range("N122").Value = -15786.6844

$ws = $wb.Worksheets.Item("WVR")
# Row 4
$ws.Range("H4").Value = 30751.25
$ws.Range("I4").Value = 5002
$ws.Range("J4").Value = 39334.332
$ws.Range("K4").Value = 5002
$ws.Range("L4").Value = 39334.332
$ws.Range("M4").Value = -4889
$ws.Range("N4").Value = -39560.332

# Row 96
$ws.Range("H96").Value = 2722.5
$ws.Range("I96").Value = 2445
$ws.Range("J96").Value = 3000
$ws.Range("K96").Value = 2445
$ws.Range("L96").Value = 3000
$ws.Range("M96").Value = -1072
$ws.Range("N96").Value = -5746
